$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1), columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update CON row (row 2), columns B:E
$ws.Range("B2").Value = 281.94773347510869
$ws.Range("C2").Value = 257.38178829125468
$ws.Range("D2").Value = 283.31402611354474
$ws.Range("E2").Value = 253.51132356186886

# Update STR row (row 3), columns B:E
$ws.Range("B3").Value = 290.43541868516547
$ws.Range("C3").Value = 249.3079991300373
$ws.Range("D3").Value = 301.305091687095
$ws.Range("E3").Value = 247.47747324598532

# Match the resulting selection reported in the saved workbook
$ws.Range("B1:E3").Select()
